$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Convertidor")

# Mark the "1/2"" rebar option as selected (aligerada de 20 completado)
# in both quantity tables on the sheet.
$ws.Range("G5").Value = 1
$ws.Range("G15").Value = 1

$excel.CalculateFullRebuild()
$wb.Save()
